# Add two new worksheets - "CaseDetailStat" and "CaseDetailStat_Message" -
# to the end of the workbook, mirroring the existing
# CypherOutput/CypherOutput_Message and StatOutput/StatOutput_Message pairs,
# but for a per-case file listing (file search for case NCATS-COP01CCB010072).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 of 2: CaseDetailStat  (tabular cypher output, like CypherOutput)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsData = $wb.Worksheets.Add($null, $lastSheet)
$wsData.Name = "CaseDetailStat"

$wsData.Range("A1").Value = "File Name"
$wsData.Range("B1").Value = "File Type"
$wsData.Range("C1").Value = "Association"
$wsData.Range("D1").Value = "Description"
$wsData.Range("E1").Value = "Format"
$wsData.Range("F1").Value = "Size"

$wsData.Range("A2").Value = "CCB010072.pdf"
$wsData.Range("B2").Value = "Pathology Report"
$wsData.Range("C2").Value = "diagnosis"
$wsData.Range("D2").Value = ""
$wsData.Range("E2").Value = "pdf"
$wsData.Range("F2").NumberFormat = "@"
$wsData.Range("F2").Value = 57.732421875
$wsData.Range("F2").ClearFormats()

# ---------------------------------------------------------------------------
# Sheet 2 of 2: CaseDetailStat_Message (connection/query log, like
# CypherOutput_Message / StatOutput_Message)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMsg = $wb.Worksheets.Add($null, $lastSheet)
$wsMsg.Name = "CaseDetailStat_Message"

$neo4jUrl = "bolt://ncias-q2251-c.nci.nih.gov:7687"
$userName = "neo4j"
$pwd = "icdcDBneo4j0"
$outputPath = "C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC05_Canine_Filter_Breed-Beagle_Neo4jData.xlsx"

$cypherCaseList = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN ['Beagle'] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(s.clinical_study_designation,'') AS ``Study Code`` , coalesce(s.clinical_study_type,'') AS  ``Study Type``, coalesce(demo.breed,'') AS Breed , coalesce(diag.disease_term,'') AS Diagnosis , coalesce(diag.stage_of_disease,'') AS ``Stage of Disease`` ,  coalesce(demo.patient_age_at_enrollment,'') AS Age , coalesce(demo.sex,'') AS Sex , coalesce(demo.neutered_indicator,'') AS  ``Neutered Status``"

$cypherBeagleCounts = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Beagle']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

$cypherFileSearch = "MATCH (f:file)-[*]->(c:case) WITH DISTINCT(f) AS f, c MATCH (f)-->(parent) WHERE c.case_id IN ['NCATS-COP01CCB010072'] RETURN f.file_name AS ``File Name`` ,f.file_type AS ``File Type``,head(labels(parent)) AS ``Association``, f.file_description AS ``Description``,f.file_format AS Format,((f.file_size)/1024) AS Size"

$cyphers = @($cypherCaseList, $cypherBeagleCounts, $cypherFileSearch)

$row = 1
foreach ($cypher in $cyphers) {
    $wsMsg.Cells.Item($row, 1).Value = "Neo4j_URL:"
    $row++
    $wsMsg.Cells.Item($row, 1).Value = $neo4jUrl
    $row++
    $wsMsg.Cells.Item($row, 1).Value = "User_name:"
    $row++
    $wsMsg.Cells.Item($row, 1).Value = $userName
    $row++
    $wsMsg.Cells.Item($row, 1).Value = "PWD:"
    $row++
    $wsMsg.Cells.Item($row, 1).Value = $pwd
    $row++
    $wsMsg.Cells.Item($row, 1).Value = "Cypher:"
    $row++
    $wsMsg.Cells.Item($row, 1).Value = $cypher
    $row++
    $wsMsg.Cells.Item($row, 1).Value = "Output:"
    $row++
    $wsMsg.Cells.Item($row, 1).Value = $outputPath
    $row++
}
